$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.621.80'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '2.524.93'
$ws.Range("E3").Value = '  +0.40%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''316.30'
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").Value = '''96.17'
$ws.Range("E6").Value = '  +0.48%  '

$ws.Range("D7").Value = '''0.575'
$ws.Range("E7").Value = '  -0.61%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = '''0.531'
$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("D10").Value = '''35.56'
$ws.Range("E10").Value = '  -1.46%  '

$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").Value = '''7.50'
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("E13").Value = '  -2.95%  '

$ws.Range("D14").Value = '2.913.07'
$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("D15").Value = '2.497.14'
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("D16").Value = '''15.06'
$ws.Range("E16").Value = '  -3.02%  '

$ws.Range("D17").Value = '''0.848'
$ws.Range("E17").Value = '  -1.06%  '

$ws.Range("D18").Value = '42.719.37'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").Value = '''6.79'

$ws.Range("D20").Value = '''12.73'
$ws.Range("E20").Value = '  -3.31%  '

$ws.Range("D21").Value = '0.0₃0960'
$ws.Range("E21").Value = '  -1.11%  '

$ws.Range("D22").Value = '''69.53'
$ws.Range("E22").Value = '  -2.54%  '

$ws.Range("D23").Value = '''250.43'
$ws.Range("E23").Value = '  -0.73%  '

$ws.Range("D24").Value = '''2.93'
$ws.Range("E24").Value = '  -2.06%  '

$ws.Range("D25").Value = '''2.04'
$ws.Range("E25").Value = '  +0.10%  '

$ws.Range("D26").Value = '''26.35'
$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("E28").Value = '  +1.80%  '

$ws.Range("D29").Value = '''40.65'
$ws.Range("E29").Value = '  +4.58%  '

$ws.Range("D30").Value = '''10.36'
$ws.Range("E30").Value = '  +2.74%  '

$ws.Range("D31").Value = '''5.90'
$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("D32").Value = '''157.27'
$ws.Range("E32").Value = '  +0.76%  '

$ws.Range("D33").Value = '''2.14'
$ws.Range("E33").Value = '  +2.59%  '

$ws.Range("D34").Value = '''2.71'
$ws.Range("E34").Value = '  +3.96%  '

$ws.Range("D35").Value = '''3.33'
$ws.Range("E35").Value = '  -0.80%  '

$ws.Range("D36").Value = '''18.88'
$ws.Range("E36").Value = '  -2.85%  '

$ws.Range("D37").Value = '''0.0781'
$ws.Range("E37").Value = '  -0.71%  '

$ws.Range("D38").Value = '''0.112'
$ws.Range("E38").Value = '  -0.67%  '

$ws.Range("E39").Value = '  -1.15%  '

$ws.Range("D40").Value = '''2.31'
$ws.Range("E40").Value = '  +8.79%  '

$ws.Range("D41").Value = '''22.43'
$ws.Range("E41").Value = '  -6.50%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0304'
$ws.Range("E42").Value = '  +1.14%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''3.81'
$ws.Range("E43").Value = '  -1.40%  '

$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").Value = '2.027.25'
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("D46").Value = '''3.25'
$ws.Range("E46").Value = '  -3.94%  '

$ws.Range("D47").Value = '''9.01'
$ws.Range("E47").Value = '  +2.05%  '

$ws.Range("D48").Value = '''84.15'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '''105.82'
$ws.Range("E49").Value = '  +3.85%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").Value = '''74.93'
$ws.Range("E50").Value = '  +2.27%  '

$ws.Range("D51").Value = '2.770.59'
$ws.Range("E51").Value = '  +0.53%  '

